# Rename "trt transposed" sheet to "trt_transposed"
$wb = $excel.ActiveWorkbook
$trtTransposed = $wb.Worksheets.Item("trt transposed")
$trtTransposed.Name = "trt_transposed"

# Make it the active sheet with a single-cell selection at A2
# (previously the selection there spanned A2:I9)
$trtTransposed.Activate() | Out-Null
$trtTransposed.Range("A2").Select() | Out-Null
